# Updated pipe thermal study - now has real thermal study things in
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 31: LH2_Cp value/unit corrected ---
$ws.Range("C31").Value = 12140
$ws.Range("D31").Value = "J/kg.K"

# --- Row 32: LH2_max_T description clarified + value corrected ---
$ws.Range("B32").Value = "maximum allowed temperature of LH2 at 21psi at the engine"
$ws.Range("C32").Value = 23

# --- Row 34: m_dot (mass flow rate of LH2 at cruise) value corrected ---
$ws.Range("C34").Value = 0.04

# --- New row 39: L_H (latent heat of vaporisation of hydrogen) ---
$ws.Range("A39").Value = "L_H"
$ws.Range("B39").Value = "latent heat of vaporisation of hydrogen"
$ws.Range("C39").Value = 449.36
$ws.Range("D39").Value = "J/mol"
$ws.Range("H39").Value = $true

# --- New row 40: mol_H (molar mass of hydrogen) ---
$ws.Range("A40").Value = "mol_H"
$ws.Range("B40").Value = "molar mass of hydrogen"
$ws.Range("C40").Formula = "=1.01*10^-3"
$ws.Range("D40").Value = "kg/mol"
$ws.Range("H40").Value = $true

# --- New row 41: boost_eta (boost pump efficiency) ---
$ws.Range("A41").Value = "boost_eta"
$ws.Range("B41").Value = "boost pump efficiency "
$ws.Range("C41").Value = 0.78
$ws.Range("H41").Value = $true

# --- New row 42: boost_m_eta (boost pump motor efficiency) ---
$ws.Range("A42").Value = "boost_m_eta"
$ws.Range("B42").Value = "boost pump motor efficiency "
$ws.Range("C42").Value = 0.86
$ws.Range("H42").Value = $true

# --- New row 43: boost_P (pressure rise from boost pump) ---
$ws.Range("A43").Value = "boost_P"
$ws.Range("B43").Value = "pressure rise from boost pump"
$ws.Range("C43").Value = 46
$ws.Range("D43").Value = "psi"
$ws.Range("H43").Value = $true

# --- New row 44: boost_power_max (maximum power required in electric boost pump) ---
$ws.Range("A44").Value = "boost_power_max"
$ws.Range("B44").Value = "maximum power required in electric boost pump"
$ws.Range("C44").Value = 3.16
$ws.Range("D44").Value = "hp"
$ws.Range("H44").Value = $true

# --- Grow the Table1 listobject + autofilter to cover the new rows ---
$table = $ws.ListObjects.Item(1)
$table.Resize($ws.Range("A1:H44"))

# --- Keep the hidden _FilterDatabase name in sync with the table range ---
$filterName = $wb.Names.Item(1)
$filterName.RefersTo = "=Sheet1!`$A`$1:`$H`$44"

# --- Restore view state (scroll position + active cell selection) ---
[void]$ws.Range("A13").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("J25").Select()
